$d = $word.ActiveDocument

$pairs = @(
    @("287÷8=35, 7", "204÷4=51, 0"),
    @("448÷9=49, 7", "589÷2=294, 1"),
    @("674÷5=134, 4", "834÷9=92, 6"),
    @("432÷6=72, 0", "115÷4=28, 3"),
    @("540÷4=135, 0", "672÷7=96, 0"),
    @("249÷8=31, 1", "909÷8=113, 5"),
    @("882÷2=441, 0", "508÷4=127, 0"),
    @("210÷2=105, 0", "662÷3=220, 2"),
    @("745÷6=124, 1", "825÷9=91, 6"),
    @("305÷2=152, 1", "155÷7=22, 1"),
    @("123÷9=13, 6", "569÷8=71, 1"),
    @("776÷3=258, 2", "896÷7=128, 0"),
    @("903÷9=100, 3", "440÷2=220, 0"),
    @("289÷8=36, 1", "213÷9=23, 6"),
    @("670÷9=74, 4", "269÷7=38, 3"),
    @("411÷7=58, 5", "134÷9=14, 8"),
    @("540÷8=67, 4", "666÷9=74, 0"),
    @("413÷4=103, 1", "785÷2=392, 1"),
    @("836÷2=418, 0", "325÷2=162, 1"),
    @("608÷6=101, 2", "866÷8=108, 2"),
    @("155÷2=77, 1", "882÷6=147, 0"),
    @("365÷4=91, 1", "252÷7=36, 0"),
    @("815÷4=203, 3", "846÷9=94, 0"),
    @("606÷5=121, 1", "554÷3=184, 2"),
    @("949÷4=237, 1", "425÷7=60, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
